$wb = $excel.ActiveWorkbook

# --- QB sheet (sheet1): Week 16 results logged, 3 reserve QBs reshuffled ---
$wsQB = $wb.Worksheets.Item("QB")

$wsQB.Range("A2").Value = "S.Ehlinger"
$wsQB.Range("A3").Value = "C.Wentz"
$wsQB.Range("A4").Value = "J.Eason"

# User finished on the QB tab, having clicked just below the table (A5)
$wsQB.Activate()
[void]$wsQB.Range("A5").Select()
